# Update Name of Algo
# Apply updated imputed values to result_data_RandomForest.xlsx

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 12.48459999999999
$ws.Range("C3").Value = -11.03229999999999
$ws.Range("E5").Value = 12.7355
$ws.Range("C14").Value = -12.332
$ws.Range("C16").Value = -12.0677
$ws.Range("E16").Value = 12.11410000000001
$ws.Range("C21").Value = -13.15930000000001
$ws.Range("C23").Value = -12.00870000000001
$ws.Range("C25").Value = -10.90239999999999
